$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Requirements sheet
# ---------------------------------------------------------------------
$wsReq = $wb.Worksheets.Item("Requirements")
$wsReq.Activate()

$wsReq.Range("B2").Value = "Output Voltage and Current"

$wsReq.Range("B3").Value = "Weight"
$wsReq.Range("C3").Value = "Maximum 200 grams"

$wsReq.Range("B4").Value = "Size: Maxiumum 3 in width, 4 in length"
$wsReq.Range("C4").Value = "Maximum 3 in width, 4 in length"

$wsReq.Range("B5").Value = "Cost"
$wsReq.Range("C5").Value = "`$200 or less"

$wsReq.Range("B6").Value = "Minimize Power Use"
$wsReq.Range("C6").Value = "Desgin Goal: The longer we can power the board the better. Minimum of 2  hours of battery life. Would like 3-4 hours"

$wsReq.Range("B7").Select()

# ---------------------------------------------------------------------
# Decision Matrix sheet
# ---------------------------------------------------------------------
$wsDM = $wb.Worksheets.Item("Decision Matrix")
$wsDM.Activate()

$wsDM.Range("B4").Value = "Can store multiple flights worth of data"
$wsDM.Range("B5").Value = "Easy to use interface"
$wsDM.Range("B6").Value = "Unit production cost < `$200"
$wsDM.Range("B8").Value = "Produces data consistantly"
$wsDM.Range("B9").Value = "3 inches width, 4 inches length"
$wsDM.Range("A10").Value = "Accuracy "
$wsDM.Range("B10").Value = "Produces accurate data"
$wsDM.Range("A11").Value = "Efficiency"
$wsDM.Range("B11").Value = "Is effcient with the power it uses."

$wsDM.Range("B3").Select()

# ---------------------------------------------------------------------
# Verification sheet
# ---------------------------------------------------------------------
$wsVer = $wb.Worksheets.Item("Verification")
$wsVer.Activate()

$wsVer.Range("B2").Value = "Output Voltage and Current"
$wsVer.Range("B3").Value = "Weight"
$wsVer.Range("B4").Value = "Size"
$wsVer.Range("B5").Value = "Cost"
$wsVer.Range("B6").Value = "Minimize Power Usage"

# Clear the old verification-method marks for rows 2-8, then re-mark
# the (possibly new) verification column for each requirement.
$wsVer.Range("C2:F8").ClearContents()

$wsVer.Range("F2").Value = "X"
$wsVer.Range("E3").Value = "X"
$wsVer.Range("C4").Value = "X"
$wsVer.Range("C5").Value = "X"
$wsVer.Range("F6").Value = "X"

$wsVer.Range("C4").Select()
